# Add team record (Wins/Losses/Ties) columns to the data sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AB1) onto the new header cells
$ws.Range("AB1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record values for every data row (2 through 45)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 93   # AD = column 30
    $ws.Cells.Item($row, 31).Value = 69   # AE = column 31
    $ws.Cells.Item($row, 32).Value = 0    # AF = column 32
}
